$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -17.64278504200363
$ws.Range("C2").Value = 2.380234231784048
$ws.Range("D2").Value = -17.64278504200363
$ws.Range("E2").Value = -17.64278504200363
$ws.Range("F2").Value = -17.64278504200363
$ws.Range("G2").Value = -17.64278504200363
$ws.Range("H2").Value = -17.64278504200363
$ws.Range("I2").Value = -17.64278504200363
$ws.Range("J2").Value = -17.64278504200363
$ws.Range("K2").Value = -17.64278504200363
$ws.Range("B3").Value = -17.64278504200363
$ws.Range("C3").Value = -17.64278504200363
$ws.Range("D3").Value = -17.64278504200363
$ws.Range("E3").Value = -17.64278504200363
$ws.Range("F3").Value = -17.64278504200363
$ws.Range("G3").Value = -17.64278504200363
$ws.Range("H3").Value = -17.64278504200363
$ws.Range("I3").Value = 2.051669771674419
$ws.Range("J3").Value = -17.64278504200363
$ws.Range("K3").Value = -17.64278504200363
$ws.Range("B4").Value = -17.64278504200363
$ws.Range("C4").Value = 2.070119737839999
$ws.Range("D4").Value = 2.919860795289067
$ws.Range("E4").Value = -17.64278504200363
$ws.Range("F4").Value = 2.485291115187787
$ws.Range("G4").Value = -17.64278504200363
$ws.Range("H4").Value = -17.64278504200363
$ws.Range("I4").Value = -17.64278504200363
$ws.Range("J4").Value = 2.573701511267257
$ws.Range("K4").Value = -17.64278504200363
$ws.Range("B5").Value = -17.64278504200363
$ws.Range("C5").Value = 0.8469571014517593
$ws.Range("D5").Value = -17.64278504200363
$ws.Range("E5").Value = -17.64278504200363
$ws.Range("F5").Value = -17.64278504200363
$ws.Range("G5").Value = 1.99715795936654
$ws.Range("H5").Value = -17.64278504200363
$ws.Range("I5").Value = -17.64278504200363
$ws.Range("J5").Value = -17.64278504200363
$ws.Range("K5").Value = -17.64278504200363
$ws.Range("B6").Value = -17.64278504200363
$ws.Range("C6").Value = -17.64278504200363
$ws.Range("D6").Value = -17.64278504200363
$ws.Range("E6").Value = -17.64278504200363
$ws.Range("F6").Value = -17.64278504200363
$ws.Range("G6").Value = -17.64278504200363
$ws.Range("H6").Value = -17.64278504200363
$ws.Range("I6").Value = -17.64278504200363
$ws.Range("J6").Value = -17.64278504200363
$ws.Range("K6").Value = -17.64278504200363
$ws.Range("B7").Value = 2.978477194230597
$ws.Range("C7").Value = -17.64278504200363
$ws.Range("D7").Value = -17.64278504200363
$ws.Range("E7").Value = -17.64278504200363
$ws.Range("F7").Value = -17.64278504200363
$ws.Range("G7").Value = -17.64278504200363
$ws.Range("H7").Value = -17.64278504200363
$ws.Range("I7").Value = -17.64278504200363
$ws.Range("J7").Value = -17.64278504200363
$ws.Range("K7").Value = -17.64278504200363
$ws.Range("B8").Value = -17.64278504200363
$ws.Range("C8").Value = -17.64278504200363
$ws.Range("D8").Value = -17.64278504200363
$ws.Range("E8").Value = 2.867520369712988
$ws.Range("F8").Value = -17.64278504200363
$ws.Range("G8").Value = -17.64278504200363
$ws.Range("H8").Value = -17.64278504200363
$ws.Range("I8").Value = -17.64278504200363
$ws.Range("J8").Value = -17.64278504200363
$ws.Range("K8").Value = -17.64278504200363
$ws.Range("B9").Value = 3.599124229664445
$ws.Range("C9").Value = -17.64278504200363
$ws.Range("D9").Value = -17.64278504200363
$ws.Range("E9").Value = -17.64278504200363
$ws.Range("F9").Value = -17.64278504200363
$ws.Range("G9").Value = -17.64278504200363
$ws.Range("H9").Value = -17.64278504200363
$ws.Range("I9").Value = -17.64278504200363
$ws.Range("J9").Value = -17.64278504200363
$ws.Range("K9").Value = -17.64278504200363
$ws.Range("B10").Value = -17.64278504200363
$ws.Range("C10").Value = -17.64278504200363
$ws.Range("D10").Value = -17.64278504200363
$ws.Range("E10").Value = -17.64278504200363
$ws.Range("F10").Value = -17.64278504200363
$ws.Range("G10").Value = -17.64278504200363
$ws.Range("H10").Value = -17.64278504200363
$ws.Range("I10").Value = 1.691027636745224
$ws.Range("J10").Value = -17.64278504200363
$ws.Range("K10").Value = 2.058692337241
$ws.Range("B11").Value = -17.64278504200363
$ws.Range("C11").Value = -17.64278504200363
$ws.Range("D11").Value = -17.64278504200363
$ws.Range("E11").Value = 1.947211655175155
$ws.Range("F11").Value = -17.64278504200363
$ws.Range("G11").Value = 2.284450384092072
$ws.Range("H11").Value = -17.64278504200363
$ws.Range("I11").Value = -17.64278504200363
$ws.Range("J11").Value = -17.64278504200363
$ws.Range("K11").Value = 1.076018715296961
$ws.Range("B12").Value = -17.64278504200363
$ws.Range("C12").Value = -17.64278504200363
$ws.Range("D12").Value = -17.64278504200363
$ws.Range("E12").Value = -17.64278504200363
$ws.Range("F12").Value = -17.64278504200363
$ws.Range("G12").Value = -17.64278504200363
$ws.Range("H12").Value = -17.64278504200363
$ws.Range("I12").Value = -17.64278504200363
$ws.Range("J12").Value = -17.64278504200363
$ws.Range("K12").Value = -17.64278504200363
$ws.Range("B13").Value = -17.64278504200363
$ws.Range("C13").Value = -17.64278504200363
$ws.Range("D13").Value = -17.64278504200363
$ws.Range("E13").Value = 1.619688362017591
$ws.Range("F13").Value = -17.64278504200363
$ws.Range("G13").Value = -17.64278504200363
$ws.Range("H13").Value = -17.64278504200363
$ws.Range("I13").Value = -17.64278504200363
$ws.Range("J13").Value = 2.268356748859713
$ws.Range("K13").Value = 1.485544148576743
$ws.Range("B14").Value = -17.64278504200363
$ws.Range("C14").Value = -17.64278504200363
$ws.Range("D14").Value = 1.586223835009056
$ws.Range("E14").Value = -17.64278504200363
$ws.Range("F14").Value = -17.64278504200363
$ws.Range("G14").Value = -17.64278504200363
$ws.Range("H14").Value = -17.64278504200363
$ws.Range("I14").Value = -17.64278504200363
$ws.Range("J14").Value = -17.64278504200363
$ws.Range("K14").Value = 2.091449912001325
$ws.Range("B15").Value = -17.64278504200363
$ws.Range("C15").Value = -17.64278504200363
$ws.Range("D15").Value = -0.3514614314439158
$ws.Range("E15").Value = -17.64278504200363
$ws.Range("F15").Value = -17.64278504200363
$ws.Range("G15").Value = -17.64278504200363
$ws.Range("H15").Value = -17.64278504200363
$ws.Range("I15").Value = -17.64278504200363
$ws.Range("J15").Value = -17.64278504200363
$ws.Range("K15").Value = -17.64278504200363
$ws.Range("B16").Value = -17.64278504200363
$ws.Range("C16").Value = -17.64278504200363
$ws.Range("D16").Value = -17.64278504200363
$ws.Range("E16").Value = -17.64278504200363
$ws.Range("F16").Value = -17.64278504200363
$ws.Range("G16").Value = -17.64278504200363
$ws.Range("H16").Value = -17.64278504200363
$ws.Range("I16").Value = -17.64278504200363
$ws.Range("J16").Value = 2.318485174833701
$ws.Range("K16").Value = -17.64278504200363
$ws.Range("B17").Value = -17.64278504200363
$ws.Range("C17").Value = 0.6243170903052194
$ws.Range("D17").Value = -0.2348783388274813
$ws.Range("E17").Value = -17.64278504200363
$ws.Range("F17").Value = -17.64278504200363
$ws.Range("G17").Value = -17.64278504200363
$ws.Range("H17").Value = 4.32192139771236
$ws.Range("I17").Value = 0.2267666320559004
$ws.Range("J17").Value = 0.9563614263683977
$ws.Range("K17").Value = -17.64278504200363
$ws.Range("B18").Value = -17.64278504200363
$ws.Range("C18").Value = -17.64278504200363
$ws.Range("D18").Value = -17.64278504200363
$ws.Range("E18").Value = -17.64278504200363
$ws.Range("F18").Value = -17.64278504200363
$ws.Range("G18").Value = -17.64278504200363
$ws.Range("H18").Value = -17.64278504200363
$ws.Range("I18").Value = 0.3111215095462049
$ws.Range("J18").Value = 1.201850238287105
$ws.Range("K18").Value = -17.64278504200363
$ws.Range("B19").Value = -17.64278504200363
$ws.Range("C19").Value = -17.64278504200363
$ws.Range("D19").Value = 1.622255350078948
$ws.Range("E19").Value = -17.64278504200363
$ws.Range("F19").Value = -17.64278504200363
$ws.Range("G19").Value = -17.64278504200363
$ws.Range("H19").Value = -17.64278504200363
$ws.Range("I19").Value = 2.292473546265237
$ws.Range("J19").Value = -17.64278504200363
$ws.Range("K19").Value = -17.64278504200363
$ws.Range("B20").Value = -17.64278504200363
$ws.Range("C20").Value = 1.85274407096351
$ws.Range("D20").Value = 2.238013641302235
$ws.Range("E20").Value = -17.64278504200363
$ws.Range("F20").Value = 3.848040868233051
$ws.Range("G20").Value = -17.64278504200363
$ws.Range("H20").Value = -17.64278504200363
$ws.Range("I20").Value = 2.410171659261255
$ws.Range("J20").Value = -17.64278504200363
$ws.Range("K20").Value = 2.736310886416056
$ws.Range("B21").Value = -17.64278504200363
$ws.Range("C21").Value = 1.864910895676352
$ws.Range("D21").Value = -17.64278504200363
$ws.Range("E21").Value = 2.529184568416931
$ws.Range("F21").Value = -17.64278504200363
$ws.Range("G21").Value = 3.47715976320588
$ws.Range("H21").Value = -17.64278504200363
$ws.Range("I21").Value = -17.64278504200363
$ws.Range("J21").Value = -17.64278504200363
$ws.Range("K21").Value = -17.64278504200363
